$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A97").Value = "2024-09-25T18:06:40Z"
$ws.Range("B97").Value = "temperature"
$ws.Range("C97").Value = "25"
$ws.Range("D97").Value = "N/A"
$ws.Range("E97").Value = "N/A"
$ws.Range("F97").Value = "N/A"
